$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "LF" column (D) entirely
$ws.Range("D1:D4").Delete()

# Remove the "LF Lag" row (row 4) entirely
$ws.Range("A4:C4").Delete()

# Update the remaining regression values
$ws.Range("B2").Value = "-0.68***"
$ws.Range("B3").Value = "9.556***"
$ws.Range("C2").Value = "-0.032***"
$ws.Range("C3").Value = "0.426***"
